$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Cells whose new value would be auto-parsed as a number by Excel; force text format first ---
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '562.27'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.67'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.349'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '25.85'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.31'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '323.64'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '66.64'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '549.55'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.15'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '153.38'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.42'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.50'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.79'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.991'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '146.51'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '19.82'

# --- Remaining cells: plain text values (not number-parseable) ---
$ws.Range("D2").Value = '61.952.02'
$ws.Range("E2").Value = '  -0.28%  '
$ws.Range("D3").Value = '2.417.72'
$ws.Range("E3").Value = '  +0.01%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("E5").Value = '  +0.69%  '
$ws.Range("E6").Value = '  -0.41%  '
$ws.Range("E8").Value = '  -0.33%  '
$ws.Range("E9").Value = '  -0.17%  '
$ws.Range("E10").Value = '  -1.18%  '
$ws.Range("E12").Value = '  -1.09%  '
$ws.Range("E13").Value = '  -1.32%  '
$ws.Range("E14").Value = '  -2.01%  '
$ws.Range("D15").Value = '2.852.68'
$ws.Range("E15").Value = '  -0.30%  '
$ws.Range("D16").Value = '61.841.16'
$ws.Range("E16").Value = '  -0.26%  '
$ws.Range("D17").Value = '2.424.45'
$ws.Range("E17").Value = '  +0.28%  '
$ws.Range("E18").Value = '  +1.28%  '
$ws.Range("E19").Value = '  -0.25%  '
$ws.Range("E20").Value = '  -1.80%  '
$ws.Range("E21").Value = '  +0.80%  '
$ws.Range("E22").Value = '  -0.06%  '
$ws.Range("E23").Value = '  +1.84%  '
$ws.Range("E24").Value = '  -0.28%  '
$ws.Range("E25").Value = '  -3.81%  '
$ws.Range("E26").Value = '  -7.49%  '
$ws.Range("D27").Value = '2.536.33'
$ws.Range("E27").Value = '  +0.61%  '
$ws.Range("E28").Value = '  -0.06%  '
$ws.Range("D29").Value = '0.0₃0929'
$ws.Range("E29").Value = '  -1.04%  '
$ws.Range("E30").Value = '  -1.78%  '
$ws.Range("E31").Value = '  -4.39%  '
$ws.Range("E32").Value = '  -0.92%  '
$ws.Range("E33").Value = '  -0.77%  '
$ws.Range("E34").Value = '  -4.18%  '
$ws.Range("E35").Value = '  -0.05%  '
$ws.Range("E36").Value = '  -1.16%  '
$ws.Range("E37").Value = '  -1.58%  '
$ws.Range("E38").Value = '  +1.77%  '
$ws.Range("E39").Value = '  -5.00%  '
$ws.Range("E40").Value = '  -1.02%  '
$ws.Range("E41").Value = '  -1.68%  '
$ws.Range("E42").Value = '  -0.84%  '
$ws.Range("E43").Value = '  -3.18%  '
$ws.Range("E44").Value = '  -5.56%  '
$ws.Range("E45").Value = '  -0.87%  '
$ws.Range("E46").Value = '  -2.64%  '
$ws.Range("E47").Value = '  -2.18%  '
$ws.Range("E48").Value = '  -0.20%  '
$ws.Range("E49").Value = '  -0.45%  '
$ws.Range("E50").Value = '  -1.29%  '
$ws.Range("E51").Value = '  +0.65%  '
